{"js": "// Updates the date line at the top of the document and every answer cell in\n// the 20x5 table of arithmetic problems, matching the target revision.\n\nconst DATE_OLD = \"2024-04-20 Saturday\";\nconst DATE_NEW = \"2024-04-21 Sunday\";\n\n// New value for every table cell, in row-major order (20 rows x 5 columns),\n// matching the table's existing row/column layout.\nconst GRID = [\n  [\"54+39=93\", \"25-7=18\", \"70-47=23\", \"76-8=68\", \"8+33=41\"],\n  [\"8+74=82\", \"46+45=91\", \"80-58=22\", \"51-43=8\", \"94-36=58\"],\n  [\"66-48=18\", \"52+39=91\", \"9+35=44\", \"83-46=37\", \"55-37=18\"],\n  [\"58+27=85\", \"48+5=53\", \"56+6=62\", \"81-69=12\", \"27+29=56\"],\n  [\"90-34=56\", \"51-25=26\", \"25-8=17\", \"7+65=72\", \"37+9=46\"],\n  [\"17-9=8\", \"52-18=34\", \"8+37=45\", \"68+28=96\", \"60-49=11\"],\n  [\"95-36=59\", \"16+28=44\", \"92-28=64\", \"29+59=88\", \"26+28=54\"],\n  [\"61-43=18\", \"70-68=2\", \"28+38=66\", \"19+56=75\", \"74-47=27\"],\n  [\"57+27=84\", \"54+28=82\", \"24+59=83\", \"6+16=22\", \"79+3=82\"],\n  [\"14+37=51\", \"56-17=39\", \"68-19=49\", \"45+37=82\", \"8+57=65\"],\n  [\"47+39=86\", \"91-63=28\", \"87-68=19\", \"41-37=4\", \"77+5=82\"],\n  [\"55+36=91\", \"67+19=86\", \"46+49=95\", \"55-47=8\", \"8+27=35\"],\n  [\"84-26=58\", \"76-19=57\", \"39+7=46\", \"87-49=38\", \"8+68=76\"],\n  [\"86-28=58\", \"27+35=62\", \"83-54=29\", \"53-6=47\", \"63-14=49\"],\n  [\"58+9=67\", \"45+46=91\", \"94-36=58\", \"17+35=52\", \"63-54=9\"],\n  [\"56-49=7\", \"6+27=33\", \"6+59=65\", \"43+38=81\", \"43+29=72\"],\n  [\"30-16=14\", \"16+59=75\", \"5+79=84\", \"73-57=16\", \"90-76=14\"],\n  [\"39+42=81\", \"52-16=36\", \"39+3=42\", \"19+65=84\", \"36+45=81\"],\n  [\"9+26=35\", \"23-16=7\", \"7+58=65\", \"42-14=28\", \"53-29=24\"],\n  [\"34+28=62\", \"50-27=23\", \"95-68=27\", \"92-36=56\", \"17+57=74\"],\n];\n\nconst body = context.document.body;\n\n// Update the date paragraph (the first paragraph in the document body).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\nif (dateParagraph.text.trim() === DATE_OLD) {\n  dateParagraph.insertText(DATE_NEW, \"Replace\");\n} else {\n  // Fall back to a direct replace-all in case paragraph indexing differs.\n  const found = body.search(DATE_OLD, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n  for (const range of found.items) {\n    range.insertText(DATE_NEW, \"Replace\");\n  }\n}\n\n// Update every cell of the (single) table with its new computed value.\nconst table = body.tables.getFirstOrNullObject();\nawait context.sync();\n\nif (!table.isNullObject) {\n  for (let r = 0; r < GRID.length; r++) {\n    for (let c = 0; c < GRID[r].length; c++) {\n      table.getCell(r, c).value = GRID[r][c];\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Updates the date line at the top of the document and every answer cell in\n# the 20x5 table of arithmetic problems, matching the target revision.\n\n$d = $word.ActiveDocument\n\n$DateOld = \"2024-04-20 Saturday\"\n$DateNew = \"2024-04-21 Sunday\"\n\n# New value for every table cell, in row-major order (20 rows x 5 columns),\n# matching the table's existing row/column layout.\n$GRID = @(\n    @(\"54+39=93\", \"25-7=18\", \"70-47=23\", \"76-8=68\", \"8+33=41\"),\n    @(\"8+74=82\", \"46+45=91\", \"80-58=22\", \"51-43=8\", \"94-36=58\"),\n    @(\"66-48=18\", \"52+39=91\", \"9+35=44\", \"83-46=37\", \"55-37=18\"),\n    @(\"58+27=85\", \"48+5=53\", \"56+6=62\", \"81-69=12\", \"27+29=56\"),\n    @(\"90-34=56\", \"51-25=26\", \"25-8=17\", \"7+65=72\", \"37+9=46\"),\n    @(\"17-9=8\", \"52-18=34\", \"8+37=45\", \"68+28=96\", \"60-49=11\"),\n    @(\"95-36=59\", \"16+28=44\", \"92-28=64\", \"29+59=88\", \"26+28=54\"),\n    @(\"61-43=18\", \"70-68=2\", \"28+38=66\", \"19+56=75\", \"74-47=27\"),\n    @(\"57+27=84\", \"54+28=82\", \"24+59=83\", \"6+16=22\", \"79+3=82\"),\n    @(\"14+37=51\", \"56-17=39\", \"68-19=49\", \"45+37=82\", \"8+57=65\"),\n    @(\"47+39=86\", \"91-63=28\", \"87-68=19\", \"41-37=4\", \"77+5=82\"),\n    @(\"55+36=91\", \"67+19=86\", \"46+49=95\", \"55-47=8\", \"8+27=35\"),\n    @(\"84-26=58\", \"76-19=57\", \"39+7=46\", \"87-49=38\", \"8+68=76\"),\n    @(\"86-28=58\", \"27+35=62\", \"83-54=29\", \"53-6=47\", \"63-14=49\"),\n    @(\"58+9=67\", \"45+46=91\", \"94-36=58\", \"17+35=52\", \"63-54=9\"),\n    @(\"56-49=7\", \"6+27=33\", \"6+59=65\", \"43+38=81\", \"43+29=72\"),\n    @(\"30-16=14\", \"16+59=75\", \"5+79=84\", \"73-57=16\", \"90-76=14\"),\n    @(\"39+42=81\", \"52-16=36\", \"39+3=42\", \"19+65=84\", \"36+45=81\"),\n    @(\"9+26=35\", \"23-16=7\", \"7+58=65\", \"42-14=28\", \"53-29=24\"),\n    @(\"34+28=62\", \"50-27=23\", \"95-68=27\", \"92-36=56\", \"17+57=74\")\n)\n\n# Update the date paragraph (the first paragraph in the document body).\n$dateRange = $d.Paragraphs(1).Range\nif ($dateRange.Text.Trim() -eq $DateOld) {\n    $dateRange.Text = $DateNew\n} else {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $DateOld\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $DateNew\n    $find.Execute($DateOld, $false, $false, $false, $false, $false, $true, 1, $false, $DateNew, 2)\n}\n\n# Update every answer cell in the (single) table with its new computed value.\n$t = $d.Tables(1)\nfor ($r = 0; $r -lt $GRID.Length; $r++) {\n    for ($c = 0; $c -lt $GRID[$r].Length; $c++) {\n        $t.Cell($r + 1, $c + 1).Range.Text = $GRID[$r][$c]\n    }\n}\n"}
